$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.8911531884013855
$ws.Range("D4").Value = 1.130117375751189
$ws.Range("D7").Value = 1.515034436477494
$ws.Range("D8").Value = 1.566005995027976
$ws.Range("D9").Value = 1.539139054070432
$ws.Range("D10").Value = 1.515727075580289
$ws.Range("D12").Value = 0.979524108265561
$ws.Range("D14").Value = 0.9973020558291696
$ws.Range("D15").Value = 1.382707248132886
$ws.Range("D20").Value = 1.341100679493786
$ws.Range("D21").Value = 1.34082899538656
$ws.Range("D23").Value = 0.8568585243570799

$ws.Range("E25").Value = 23.5
$ws.Range("E29").Value = 19
$ws.Range("E31").Value = 12
$ws.Range("E34").Value = 26
